$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New "GDP per Capita" Data-column values for years 1952-2016 (row 2 = year 1952 ... row 66 = year 2016).
# These are stored as text in the workbook (shared strings), so a leading apostrophe is used to stop
# Excel from auto-converting the numeric-looking text into a number.
$values = @(
    "2032", "2348", "2236", "2584", "2447", "3147", "3022", "3885", "3948", "4038",
    "4237", "4713", "5166", "5231", "5649", "5703", "5761", "6379", "6515", "7084",
    "7347", "7716", "8322", "8498", "8877", "9529", "9958", "10772", "10724", "10888",
    "10949", "10707", "10866", "10839", "11180", "10973", "10793", "10963", "9754",
    "8858.00687718841", "6567.66979131886", "4673.64158257596", "4910.81639358195",
    "5325.48823239337", "5710.196717778", "6288.02412151558", "6603.71117893823",
    "5532.54132555933", "6059.58302217881", "6542.67954226357", "7155.57619129481",
    "7637.83224410049", "8512.16789709297", "9198.41447412235", "9902.42608751295",
    "10768.9423788207", "11650.228005412", "11576.5888058709", "11963.5550659209",
    "12463", "12384", "12764", "12606", "12765", "13184"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = "'" + $values[$i]
}

# Add new rows for years 2011-2016 (rows 61-66), same layout as existing rows.
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)
for ($j = 0; $j -lt $newYears.Length; $j++) {
    $row = 61 + $j
    $ws.Cells.Item($row, 1).Value = 688
    $ws.Cells.Item($row, 2).Value = "Serbia"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $newYears[$j]
    $ws.Cells.Item($row, 5).Value = "'" + $values[59 + $j]
}
